$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the "Chain (1ft)" line item to the more specific "04B Chain (1ft)" ---
# (keep the leading space exactly as authored) and drop its special centered
# formatting on the link cell so it reads like a normal (left-aligned) cell.
$ws.Range("A6").Value = " 04B Chain (1ft)"
$ws.Range("B6").ClearFormats()

# --- Add the new "Large Cable Transmission System" item: a Bowden Tube ---
$ws.Range("A20").Value = "Bowden Tube"
$ws.Range("B20").Value = "https://www.amazon.com/Jagwire-Sport-Housing-Slick-Lube-Titanium/dp/B085NBZMJS/ref=sr_1_1?crid=1A5WX5ADQYM0Y&dib=eyJ2IjoiMSJ9.CWV7EelBoN67bHKqG_VMGaAHwKq3lWAqkaCYopLdT43GyH4CDNeyWoQ_bFV_YrbZmhmwmsofP69GRzCCWYW_ULIkinZgZrdky8EGo_FPRa2GDLIPcrjwSu8T1nDFsZ03wyuLLatilRsdmpFkqvecV8S7AOhy1XjGVW6Ztcl1kgVL4_2zjOWbooP9z_kj4elJxdrMd7yL-uhr5ZcYm6F_Z725qnNy_c0-wNfaDt_xxhvLSabJzp2Ta9HzGpozVxtWGMo89NRT0qD_9iY1NxHHBedJp86w9TXZ0OfMm9s45WQ.7sQYhAOJzDOGkHvrSoGgv6F0CIzum6_7bjb1Zjv8zzc&dib_tag=se&keywords=Jagwire+Brake+Housing+CGX-SL+Slick-Lube+5+mm+%2810+m%29&qid=1726249759&s=sporting-goods&sprefix=jagwire+brake+housing+cgx-sl+slick-lube+5+mm+10+m+%2Csporting%2C99&sr=1-1"
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 61.97

$ws.Range("A20").HorizontalAlignment = -4108
$ws.Range("C20").HorizontalAlignment = -4108
$ws.Range("D20").HorizontalAlignment = -4108

# Move the active selection to the newly-added row, matching where the
# author left the cursor after typing the new entry.
$ws.Range("A20:D20").Select()
